$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(2).ColumnWidth = 31.453125
$w = $ws1.Columns.Item(2).ColumnWidth()
Write-Host "Column B width after set 31.453125:" $w
$w2 = $ws1.Columns.Item(2).Width()
Write-Host "Column B Width prop:" $w2
